$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.954043392504931
$ws.Range("C2").Value = 0.02019987732393777

$ws.Range("B3").Value = 0.8860453648915186
$ws.Range("C3").Value = 0.07234400031150691

$ws.Range("B4").Value = 0.9245069033530573
$ws.Range("C4").Value = 0.09671118753893371

$ws.Range("B5").Value = 0.777120315581854
$ws.Range("C5").Value = 0.09861735698224811

$ws.Range("B6").Value = 0.892455621301775
$ws.Range("C6").Value = 0.04702705532228531
